$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.648.38"
$ws.Range("E2").Value = "  -2.17%  "
$ws.Range("D3").Value = "1.794.47"
$ws.Range("E3").Value = "  -2.00%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.26"
$ws.Range("E5").Value = "  -1.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5877"
$ws.Range("E6").Value = "  -2.51%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2761"
$ws.Range("E8").Value = "  -1.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06740"
$ws.Range("E9").Value = "  -4.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.14"
$ws.Range("E10").Value = "  -1.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07530"
$ws.Range("E11").Value = "  -1.77%  "
$ws.Range("D12").Value = "1.796.57"
$ws.Range("E12").Value = "  -1.80%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.782"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6138"
$ws.Range("E14").Value = "  -2.44%  "
$ws.Range("D15").Value = "2.038.07"
$ws.Range("E15").Value = "  -1.98%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "75.21"
$ws.Range("E16").Value = "  -4.97%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000009028"
$ws.Range("E17").Value = "  -8.68%  "
$ws.Range("D18").Value = "28.630.85"
$ws.Range("E18").Value = "  -2.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.455"
$ws.Range("E19").Value = "  -6.85%  "
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "209.46"
$ws.Range("E21").Value = "  -6.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.47"
$ws.Range("E22").Value = "  -2.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.809"
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.11"
$ws.Range("E25").Value = "  -1.97%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.037"
$ws.Range("E26").Value = "  +0.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1257"
$ws.Range("E27").Value = "  -3.90%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.39"
$ws.Range("E28").Value = "  -1.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.420"
$ws.Range("E29").Value = "  -4.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06099"
$ws.Range("E30").Value = "  -4.27%  "
$ws.Range("E31").Value = "  -1.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.807"
$ws.Range("E32").Value = "  +0.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.779"
$ws.Range("E33").Value = "  -1.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.729"
$ws.Range("E34").Value = "  -0.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.046"
$ws.Range("E35").Value = "  -5.86%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6391"
$ws.Range("E36").Value = "  -1.32%  "
$ws.Range("E37").Value = "  -1.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.712"
$ws.Range("E38").Value = "  -1.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.404"
$ws.Range("E39").Value = "  -2.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01694"
$ws.Range("E40").Value = "  -3.39%  "
$ws.Range("D41").Value = "1.141.73"
$ws.Range("E41").Value = "  -6.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8787"
$ws.Range("E42").Value = "  -2.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.06"
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("D45").Value = "1.945.67"
$ws.Range("E45").Value = "  -2.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "59.84"
$ws.Range("E46").Value = "  -4.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000110"
$ws.Range("E47").Value = "  -5.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.581"
$ws.Range("E48").Value = "  +0.21%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.366"
$ws.Range("E49").Value = "  -2.60%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05482"
$ws.Range("E50").Value = "  -0.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4482"
$ws.Range("E51").Value = "  -1.66%  "
